$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.087.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.429.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '411.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.755'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000202'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +57.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.432.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +14.63%  '
$ws.Range("E18").Value = '  +4.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.143.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '402.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +27.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '90.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.97%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '33.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.53%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '44.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.55%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.172'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.41%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0503'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '52.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.62%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.10%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.315'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.120.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0371'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.21%  '
